$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C5").Value = "Please refer Working template 1"
$ws.Range("C28").Value = "Please refer Working template 1"

$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollColumn = 2
